# edit.ps1
# Re-themes the "IT_KPI_Dashboard" template from an Information-Technology /
# Cloud-Infrastructure narrative to an Artificial-Intelligence & Machine-Learning
# narrative, per the authoritative diff.
#
# Touches two worksheets:
#   1. "Instructions & User Guide"  - title + KPI glossary rows
#   2. "KPI Dashboard"              - header/title rows + KPI table (rows 8-22)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")
$ws2 = $wb.Worksheets.Item("KPI Dashboard")

# ---------------------------------------------------------------------------
# Sheet 1: "Instructions & User Guide"
# ---------------------------------------------------------------------------

$ws1.Range("A1").Value  = "Artificial Intelligence and Machine Learning KPI Dashboard - User Guide & Instructions"

# Row 2 is a blank spacer row that must be preserved (not pruned) across the
# save even though it carries no cell data - touching a benign, already-default
# row-level property keeps it "present" without materialising any cell.
$ws1.Rows.Item(2).OutlineLevel = 0

$ws1.Range("A19").Value = "Model Accuracy Rate"
$ws1.Range("B19").Value = "Key performance indicator for artificial intelligence and machine learning"

$ws1.Range("B22").Value = "Key performance indicator for artificial intelligence and machine learning"

$ws1.Range("B24").Value = "Key performance indicator for artificial intelligence and machine learning"

# ---------------------------------------------------------------------------
# Sheet 2: "KPI Dashboard"
# ---------------------------------------------------------------------------

$ws2.Range("A1").Value = "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING - KPI DASHBOARD"
$ws2.Range("A2").Value = "Project: AI/ML Implementation"

# KPI001 (row 8)
$ws2.Range("B8").Value = "Model Accuracy Rate"
$ws2.Range("I8").Value = "Chief Data Officer"
$ws2.Range("K8").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI002 (row 9)
$ws2.Range("B9").Value = "Data Quality Score"
$ws2.Range("I9").Value = "Data Scientists"
$ws2.Range("K9").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI003 (row 10)
$ws2.Range("B10").Value = "User Adoption Rate"
$ws2.Range("I10").Value = "ML Engineers"
$ws2.Range("K10").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI004 (row 11)
$ws2.Range("I11").Value = "Business Analysts"
$ws2.Range("K11").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI005 (row 12)
$ws2.Range("K12").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI006 (row 13)
$ws2.Range("K13").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI007 (row 14)
$ws2.Range("I14").Value = "Chief Data Officer"
$ws2.Range("K14").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI008 (row 15)
$ws2.Range("I15").Value = "Data Scientists"
$ws2.Range("K15").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI009 (row 16)
$ws2.Range("I16").Value = "ML Engineers"
$ws2.Range("K16").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI010 (row 17)
$ws2.Range("I17").Value = "Business Analysts"
$ws2.Range("K17").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI011 (row 18)
$ws2.Range("K18").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI012 (row 19)
$ws2.Range("K19").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI013 (row 20)
$ws2.Range("I20").Value = "Chief Data Officer"
$ws2.Range("K20").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI014 (row 21)
$ws2.Range("I21").Value = "Data Scientists"
$ws2.Range("K21").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"

# KPI015 (row 22)
$ws2.Range("I22").Value = "ML Engineers"
$ws2.Range("K22").Value = "Critical KPI for Artificial Intelligence and Machine Learning success"
